# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 3486
$ws.Range('F5').Value = 8165
$ws.Range('F10').Value = 1106
$ws.Range('F12').Value = 163
$ws.Range('F15').Value = 5
$ws.Range('F17').Value = 747
$ws.Range('F19').Value = 536
$ws.Range('F22').Value = 1377
$ws.Range('F23').Value = 6870
$ws.Range('F24').Value = 121
$ws.Range('F25').Value = 54071
$ws.Range('F26').Value = 4216
$ws.Range('F28').Value = 1016
$ws.Range('F32').Value = 854
$ws.Range('F35').Value = 2042
$ws.Range('F39').Value = 1086
$ws.Range('F40').Value = 502
$ws.Range('F42').Value = 171
$ws.Range('F44').Value = 687
$ws.Range('F48').Value = 31

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F9').Value = 146
$ws.Range('F10').Value = 25
$ws.Range('F11').Value = 45
$ws.Range('F12').Value = 42
$ws.Range('F16').Value = 165
$ws.Range('F30').Value = 75

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 2247
$ws.Range('F5').Value = 1510
$ws.Range('F9').Value = 9290
$ws.Range('F10').Value = 1582
$ws.Range('F12').Value = 65
$ws.Range('F15').Value = 115

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 2247
$ws.Range('F4').Value = 8165
$ws.Range('F6').Value = 1582
$ws.Range('F14').Value = 163
$ws.Range('F16').Value = 5
$ws.Range('F17').Value = 536
$ws.Range('F19').Value = 6870
$ws.Range('F20').Value = 121
$ws.Range('F21').Value = 54070
$ws.Range('F22').Value = 146
$ws.Range('F23').Value = 146
$ws.Range('F24').Value = 25
$ws.Range('F25').Value = 45
$ws.Range('F26').Value = 4216
$ws.Range('F28').Value = 42
$ws.Range('F33').Value = 2042
$ws.Range('F36').Value = 1087
$ws.Range('F41').Value = 687
$ws.Range('F46').Value = 75
$ws.Range('F47').Value = 31

# --- Sheet: 全部类型 (row 9 and row 10 content replace) ---
$ws = $wb.Worksheets.Item('全部类型')

# Row 9 new values
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '2024-09-21'
$ws.Range("C9").Value = '上海·第十八届Redamancy动漫游戏嘉年华'
$ws.Range("D9").Value = '古方路与沪闵路交叉口正南方向136米 PK西餐厅运动娱乐中心'
$ws.Range("E9").Value = '2024.09.21 10:00-09.22 17:00'
$ws.Range("F9").Value = 517
$ws.Range("G9").Value = 78
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=92346'
$ws.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202409/XXx4B6C71726067047122.jpeg'

# Row 10 new values
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '2024-09-22'
$ws.Range("C10").Value = '上海·创世次元同人only展X探索无界，共赴次元之旅X一场跨越次元的盛宴正等待着每一位旅行者的心跳加速（免费活动）'
$ws.Range("D10").Value = '世纪大道1217号 百联世纪购物中心'
$ws.Range("E10").Value = '2024.09.22 13:30-09.22 17:30'
$ws.Range("F10").Value = 1106
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=92201'
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202409/x5RsyeXz1725883308237.jpeg'
